# Daily automatic data refresh for the EPEX Spot prices workbook.
#
# 1) "Prix Spot" sheet: append a new day column (AK) for "20-jul" with its
#    24 hourly prices.
# 2) "Gaz" sheet: append a new row (34) for 2025-07-18.
# 3) "CO2" sheet: append a new row (34) for 2025-07-18.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot" — new column AK ("20-jul")
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Prix Spot")

# AK1 is a header cell like the rest of row 1 (bold / bordered / centered
# style). Copy the formatting from the previous header cell (AJ1) so the
# new column reuses the existing header style instead of creating a new
# one, then set its text.
$ws.Range("AJ1").Copy()
$ws.Range("AK1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("AK1").Value = "20-jul"

# Hourly values for 20-jul, rows 2 (00-01) through 25 (23-24).
$akValues = @(
    51.61, 43.01, 25.5, 22.64, 21.83, 19.99, 16.32, 20.76,
    3.6, 0.65, 0, 0, 0, -0.01, -0.01, 0, 3.08, 0,
    5.99, 11.17, 43.81, 50.8, 78.33, 74.93000000000001
)
for ($i = 0; $i -lt $akValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 37).Value = $akValues[$i]   # column 37 = AK
}

# ---------------------------------------------------------------------
# Sheet "Gaz" — new row 34 (2025-07-18)
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Writing an ISO-looking date string straight into .Value lets Excel
# auto-convert it into a date serial (and forces a new number-format
# style), but the existing "Date" column stores plain text. Build the
# text via a formula first and paste back as a literal value so the cell
# keeps the default (unstyled) text type, matching the rest of the column.
$wsGaz.Range("A34").Formula = "=""2025-07-18"""
$wsGaz.Range("A34").Copy()
$wsGaz.Range("A34").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$wsGaz.Range("B34").Value = 32.85

# ---------------------------------------------------------------------
# Sheet "CO2" — new row 34 (2025-07-18)
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A34").Formula = "=""2025-07-18"""
$wsCO2.Range("A34").Copy()
$wsCO2.Range("A34").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$wsCO2.Range("B34").Value = 69.2

Write-Output "Added 20-jul column to Prix Spot and 2025-07-18 rows to Gaz/CO2"
